# Apply the OOXML diff: update header labels (row 1) to human-readable
# Spanish text, and re-shuffle the measure/dimension metadata rows
# (rows 2-4) across columns C, D, E, F, G, H, I, J so that the DSD
# metadata lines up correctly with each column's semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers
$ws.Range("A1").Value = "Personas residentes viviendas familiares"
$ws.Range("B1").Value = "Nivel estudios (agregado)"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Nivel estudios (detalle)"
$ws.Range("E1").Value = "Comarca código"
$ws.Range("F1").Value = "Provincia código"
$ws.Range("G1").Value = "Aragón"
$ws.Range("H1").Value = "Municipio código"
$ws.Range("I1").Value = "Provincia nombre"
$ws.Range("J1").Value = "Municipio nombre"

# Row 2 - concept / property reference
$ws.Range("A2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("B2").Value = "iaest-measure:nivel-estudios-agregado"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:nivel-estudios-detalle"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "null"
$ws.Range("I2").Value = "sdmx-dimension:refArea"
$ws.Range("J2").Value = "sdmx-dimension:refArea"

# Row 3 - component type (measure/dimension)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "null"
$ws.Range("I3").Value = "dim"
$ws.Range("J3").Value = "dim"

# Row 4 - data type / codelist URI
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "null"
$ws.Range("I4").Value = "URI-Provincia"
$ws.Range("J4").Value = "URI-Municipio"
